# Generate Report for Archive
# - Update status text from "Ready for handoff" to "In Translation"
#   (affects Overview!E2:F2, zh-cn!C2 and de-de!C2, which all shared the
#   same string).
# - The status column(s) narrow to fit the new, shorter text
#   (Overview columns E:F and the "Status" column C on the zh-cn / de-de
#   sheets), mirroring the column-width change Excel performs when the
#   cell content that drove the existing width becomes shorter.

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"
$newWidth = 12.576851254417766

$ws1 = $wb.Worksheets.Item("Overview")
$ws1.Cells.Replace($oldStatus, $newStatus)
$ws1.Range("E1").EntireColumn.ColumnWidth = $newWidth
$ws1.Range("F1").EntireColumn.ColumnWidth = $newWidth

$ws2 = $wb.Worksheets.Item("zh-cn")
$ws2.Cells.Replace($oldStatus, $newStatus)
$ws2.Range("C1").EntireColumn.ColumnWidth = $newWidth

$ws3 = $wb.Worksheets.Item("de-de")
$ws3.Cells.Replace($oldStatus, $newStatus)
$ws3.Range("C1").EntireColumn.ColumnWidth = $newWidth
